$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) updates ---

# Values whose new text still looks like a genuine number (e.g. "226.89")
# would otherwise be auto-converted by Excel into a numeric cell. Force
# column D to Text format first so these stay plain text, matching the
# original inline-string cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"
$ws.Range("D5").Value = '226.89'
$ws.Range("D8").Value = '31.94'
$ws.Range("D10").Value = '0.0690'
$ws.Range("D11").Value = '0.0946'
$ws.Range("D18").Value = '68.14'
$ws.Range("D19").Value = '245.60'
$ws.Range("D24").Value = '2.06'
$ws.Range("D25").Value = '161.06'
$ws.Range("D26").Value = '7.18'
$ws.Range("D30").Value = '1.24'
$ws.Range("D36").Value = '0.649'
$ws.Range("D37").Value = '2.42'
$ws.Range("D45").Value = '0.0509'
$ws.Range("D46").Value = '6.06'
$ws.Range("D50").Value = '105.98'

# Restore the default (Normal) style on column D so no stray formatting remains
$priceRange.Style = "Normal"

# These new values already contain multiple "." separators, so Excel keeps
# them as text on their own; use Replace so the existing (default) cell
# style is left completely untouched.
$ws.Range("D2").Replace('34.159.72', '34.192.59') | Out-Null
$ws.Range("D3").Replace('1.791.02', '1.791.16') | Out-Null
$ws.Range("D12").Replace('2.049.34', '2.049.75') | Out-Null
$ws.Range("D14").Replace('1.788.34', '1.786.57') | Out-Null
$ws.Range("D15").Replace('34.110.94', '34.142.75') | Out-Null
$ws.Range("D35").Replace('1.460.69', '1.453.11') | Out-Null
$ws.Range("D49").Replace('1.950.41', '1.950.85') | Out-Null

# --- Volume(1h) (E) updates ---
# Padded with two leading/trailing spaces, so Excel always keeps these as text.
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  -1.40%  '
$ws.Range("E9").Value = '  +1.08%  '
$ws.Range("E10").Value = '  -2.19%  '
$ws.Range("E11").Value = '  +0.93%  '
$ws.Range("E12").Value = '  +0.47%  '
$ws.Range("E13").Value = '  +0.69%  '
$ws.Range("E14").Value = '  -0.27%  '
$ws.Range("E15").Value = '  +0.33%  '
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("E17").Value = '  +0.61%  '
$ws.Range("E18").Value = '  +0.31%  '
$ws.Range("E19").Value = '  +0.81%  '
$ws.Range("E20").Value = '  -0.73%  '
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("E24").Value = '  +0.29%  '
$ws.Range("E25").Value = '  +0.45%  '
$ws.Range("E26").Value = '  +0.86%  '
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("E28").Value = '  +0.50%  '
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("E30").Value = '  -0.26%  '
$ws.Range("E31").Value = '  +1.26%  '
$ws.Range("E32").Value = '  +0.18%  '
$ws.Range("E33").Value = '  +2.68%  '
$ws.Range("E34").Value = '  -0.60%  '
$ws.Range("E35").Value = '  +4.01%  '
$ws.Range("E36").Value = '  -1.83%  '
$ws.Range("E37").Value = '  +7.21%  '
$ws.Range("E38").Value = '  +2.31%  '
$ws.Range("E39").Value = '  -0.25%  '
$ws.Range("E40").Value = '  +2.96%  '
$ws.Range("E42").Value = '  +0.66%  '
$ws.Range("E43").Value = '  +0.70%  '
$ws.Range("E44").Value = '  +0.41%  '
$ws.Range("E45").Value = '  +2.62%  '
$ws.Range("E46").Value = '  +3.55%  '
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("E48").Value = '  -2.89%  '
$ws.Range("E49").Value = '  +0.71%  '
$ws.Range("E50").Value = '  -2.16%  '
$ws.Range("E51").Value = '  +0.13%  '
